$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains a weekly price log for "Espinaca" (Femacal de La Calera).
# A new weekly record is inserted right before the current row 368, pushing
# every following row (old 368..430) down by one (new 369..431).

$ws.Rows(368).Insert()

$ws.Range("A368").Value = 3
$ws.Range("B368").Value = "Femacal de La Calera"
$ws.Range("C368").Value = "Coquimbo"
$ws.Range("D368").Value = 44406
$ws.Range("E368").Value = 5
$ws.Range("F368").Value = 100112012
$ws.Range("G368").Value = "Espinaca"
$ws.Range("H368").Value = "Sin especificar"
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 100
$ws.Range("K368").Value = 1300
$ws.Range("L368").Value = 1400
$ws.Range("M368").Value = 1350
$ws.Range("N368").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O368").Value = "Provincia de Quillota"
$ws.Range("P368").Value = 450
$ws.Range("Q368").Value = 3
$ws.Range("R368").Value = "Hortaliza"
